# Generate Report for Handback
#
# This script mirrors the "handback" step of the localization pipeline:
#   * the Overview / per-locale Status cells flip from "Ready for handoff"
#     to "Handed back: in sync with en-US"
#   * each locale sheet's "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" cells get populated (they were blank /
#     placeholder before)
#   * a couple of columns get widened so the new, longer text fits
#   * the new "Latest Target File" cell becomes a hyperlink to the source
#     markdown doc, just like column A already is

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$mdFile          = "afb0c21f-1429-4d00-ae16-26a10c9a1f23.md"
$mdHyperlinkUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/12c7d64548f8094e6b05d67ee44f8866ef6bb617/e2e/afb0c21f-1429-4d00-ae16-26a10c9a1f23.md"
$statusHandedBack = "Handed back: in sync with en-US"

$hyperlinkColor = 15570276  # BGR for RGB FF6495ED - matches the workbook's existing HyperLink style

# ---------------------------------------------------------------------
# Overview sheet: Status columns (E/F) now read "Handed back..."
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack

# Widen the Status columns on the Overview sheet so the longer text fits
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $statusHandedBack

# Widen Status column (C) and the two "Latest Target/Handback File" columns
$wsZhCn.Columns.Item(3).ColumnWidth = 29.17
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

# Latest Target File (I2): becomes a hyperlink to the source .md, same display text as A2
$i2 = $wsZhCn.Range("I2")
$wsZhCn.Hyperlinks.Add($i2, $mdHyperlinkUrl, "", "", $mdFile)
$i2.Value = $mdFile
$i2.Style = "HyperLink"
$i2.Font.Underline = 2
$i2.Font.Color = $hyperlinkColor

# Latest Handback File (J2): the generated zh-cn xliff file name
$wsZhCn.Range("J2").Value = "afb0c21f-1429-4d00-ae16-26a10c9a1f23.46890682b58ddfdbeb45b52fe03d5d5b7ff4d997.zh-cn.xlf"

# Latest Handback DateTime (K2)
$wsZhCn.Range("K2").Value = "2016-08-23 15:14:40"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $statusHandedBack

# Widen Status column (C) and the two "Latest Target/Handback File" columns
$wsDeDe.Columns.Item(3).ColumnWidth = 29.17
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17

# Latest Target File (I2): becomes a hyperlink to the source .md, same display text as A2
$i2de = $wsDeDe.Range("I2")
$wsDeDe.Hyperlinks.Add($i2de, $mdHyperlinkUrl, "", "", $mdFile)
$i2de.Value = $mdFile
$i2de.Style = "HyperLink"
$i2de.Font.Underline = 2
$i2de.Font.Color = $hyperlinkColor

# Latest Handback File (J2): the generated de-de xliff file name
$wsDeDe.Range("J2").Value = "afb0c21f-1429-4d00-ae16-26a10c9a1f23.46890682b58ddfdbeb45b52fe03d5d5b7ff4d997.de-de.xlf"

# Latest Handback DateTime (K2)
$wsDeDe.Range("K2").Value = "2016-08-23 15:14:48"
